# B6-PowerPoint.pptx update (Fri, Jul 31, 2020 11:06:04 PM)
#
# 1) Re-style the three summary tables (slides 14-16) from the default
#    "Table_0" style to the "Medium Style 2 - Accent 1"-class style
#    {EE4E6D26-E5D8-4C66-B236-63630D636069}.
# 2) Swap the deck's colour scheme back to the stock Office palette
#    (the design gallery only keeps one live colour scheme in this
#    session, so repainting it with the Office values reproduces the
#    "Office Theme" <-> "Integral" theme swap seen in the saved package).

$p = $ppt.ActivePresentation

$newTableStyle = "{EE4E6D26-E5D8-4C66-B236-63630D636069}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
